# Weekly update: insert a new week's price record for Berenjena
# (Mercado Mayorista Lo Valledor de Santiago) at row 297, pushing the
# existing rows 297-338 down to 298-339.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 297 (shifts rows 297:338 -> 298:339,
# carrying formatting down from the row above, same as Excel's UI).
$ws.Rows.Item(297).Insert()

# Populate the newly inserted row with this week's data.
$ws.Cells.Item(297, 1).Value = 6
$ws.Cells.Item(297, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(297, 3).Value = "Metropolitana"
$ws.Cells.Item(297, 4).Value = 45131
$ws.Cells.Item(297, 5).Value = 13
$ws.Cells.Item(297, 6).Value = 100112001
$ws.Cells.Item(297, 7).Value = "Berenjena"
$ws.Cells.Item(297, 8).Value = "Sin especificar"
$ws.Cells.Item(297, 9).Value = "Primera"
$ws.Cells.Item(297, 10).Value = 410
$ws.Cells.Item(297, 11).Value = 7000
$ws.Cells.Item(297, 12).Value = 8000
$ws.Cells.Item(297, 13).Value = 7366
$ws.Cells.Item(297, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(297, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(297, 16).Value = 147
$ws.Cells.Item(297, 17).Value = 50
$ws.Cells.Item(297, 18).Value = "Hortaliza"
